$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 32999
$ws.Range("J26").Value = 32999
$ws.Range("L26").Value = 32999
$ws.Range("N26").Value = -33687
$ws.Range("H62").Value = 1221.25
$ws.Range("I62").Value = 961.6667
$ws.Range("K62").Value = 961.6667
$ws.Range("M62").Value = -337.6667
$ws.Range("H65").Value = 1221.25
$ws.Range("I65").Value = 961.6667
$ws.Range("K65").Value = 4808.3335
$ws.Range("M65").Value = -1688.3335
$ws.Range("H123").Value = 40792.312
$ws.Range("J123").Value = 40792.312
$ws.Range("L123").Value = 40792.312
$ws.Range("N123").Value = -50592.312
$ws.Range("H137").Value = 937248.9399999999
$ws.Range("I137").Value = 2510463.8
$ws.Range("J137").Value = 3152.7188
$ws.Range("K137").Value = 7531391.399999999
$ws.Range("L137").Value = 9458.1564
$ws.Range("M137").Value = -7528841.399999999
$ws.Range("N137").Value = -14558.1564
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1431.0294
$ws.Range("I2").Value = 1432.3334
$ws.Range("K2").Value = 1432.3334
$ws.Range("M2").Value = -1319.3334
$ws.Range("H61").Value = 1638.875
$ws.Range("I61").Value = 1638.875
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1638.875
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1426.875
$ws.Range("N61").Value = ""
$ws.Range("H74").Value = 465111.03
$ws.Range("J74").Value = 2802.4
$ws.Range("L74").Value = 2802.4
$ws.Range("N74").Value = -4550.4
$ws.Range("H77").Value = 465111.03
$ws.Range("J77").Value = 2802.4
$ws.Range("L77").Value = 14012
$ws.Range("N77").Value = -22748
$ws.Range("H116").Value = 1431.0294
$ws.Range("I116").Value = 1432.3334
$ws.Range("K116").Value = 1432.3334
$ws.Range("M116").Value = 861.6666
$ws.Range("H132").Value = 2504.7727
$ws.Range("I132").Value = 1534.2307
$ws.Range("J132").Value = 3906.6667
$ws.Range("K132").Value = 4602.6921
$ws.Range("L132").Value = 11720.0001
$ws.Range("M132").Value = -2072.6921
$ws.Range("N132").Value = -16780.0001
$ws.Range("H136").Value = 1638.875
$ws.Range("I136").Value = 1638.875
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4916.625
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2366.625
$ws.Range("N136").Value = ""
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1431.0294
$ws.Range("I3").Value = 1432.3334
$ws.Range("K3").Value = 1432.3334
$ws.Range("M3").Value = -1318.3334
$ws.Range("H21").Value = 30385.5
$ws.Range("J21").Value = 30385.5
$ws.Range("L21").Value = 30385.5
$ws.Range("N21").Value = -30857.5
$ws.Range("H134").Value = 3640.24
$ws.Range("I134").Value = 1496.6154
$ws.Range("J134").Value = 5962.5
$ws.Range("K134").Value = 4489.8462
$ws.Range("L134").Value = 17887.5
$ws.Range("M134").Value = -1954.8462
$ws.Range("N134").Value = -22957.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2523.796
$ws.Range("I31").Value = 1023.619
$ws.Range("J31").Value = 3648.9285
$ws.Range("K31").Value = 1023.619
$ws.Range("L31").Value = 3648.9285
$ws.Range("M31").Value = -728.619
$ws.Range("N31").Value = -4238.9285
$ws.Range("H34").Value = 2523.796
$ws.Range("I34").Value = 1023.619
$ws.Range("J34").Value = 3648.9285
$ws.Range("K34").Value = 1023.619
$ws.Range("L34").Value = 3648.9285
$ws.Range("M34").Value = -821.619
$ws.Range("N34").Value = -4052.9285
$ws.Range("H58").Value = 2568.4849
$ws.Range("I58").Value = 1481.3478
$ws.Range("K58").Value = 1481.3478
$ws.Range("M58").Value = -1278.3478
$ws.Range("H134").Value = 2511.7273
$ws.Range("I134").Value = 1181
$ws.Range("K134").Value = 3543
$ws.Range("M134").Value = -1008
$ws.Range("H136").Value = 2568.4849
$ws.Range("I136").Value = 1481.3478
$ws.Range("K136").Value = 4444.0434
$ws.Range("M136").Value = -1894.0434
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30228
$ws.Range("H68").Value = 1141.381
$ws.Range("I68").Value = 725.4286
$ws.Range("J68").Value = 1973.2858
$ws.Range("K68").Value = 2176.2858
$ws.Range("L68").Value = 5919.857400000001
$ws.Range("M68").Value = -1365.2858
$ws.Range("N68").Value = -7541.857400000001
$ws.Range("H71").Value = 1141.381
$ws.Range("I71").Value = 725.4286
$ws.Range("J71").Value = 1973.2858
$ws.Range("K71").Value = 6528.8574
$ws.Range("L71").Value = 17759.5722
$ws.Range("M71").Value = -2472.8574
$ws.Range("N71").Value = -25871.5722
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2925.5
$ws.Range("I102").Value = 2322.7646
$ws.Range("J102").Value = 3857
$ws.Range("K102").Value = 2322.7646
$ws.Range("L102").Value = 3857
$ws.Range("M102").Value = -700.7646
$ws.Range("N102").Value = -7101
$ws.Range("H122").Value = 6138.048
$ws.Range("I122").Value = 2990.9092
$ws.Range("J122").Value = 9599.9
$ws.Range("K122").Value = 8972.7276
$ws.Range("L122").Value = 28799.7
$ws.Range("M122").Value = -6522.7276
$ws.Range("N122").Value = -33699.7
$ws.Range("H126").Value = 3956.2917
$ws.Range("I126").Value = 2938.2942
$ws.Range("K126").Value = 8814.882599999999
$ws.Range("M126").Value = -6344.882599999999
$ws.Range("H132").Value = 3954.5356
$ws.Range("I132").Value = 2953.5625
$ws.Range("J132").Value = 5289.1665
$ws.Range("K132").Value = 8860.6875
$ws.Range("L132").Value = 15867.4995
$ws.Range("M132").Value = -6330.6875
$ws.Range("N132").Value = -20927.4995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3661.7693
$ws.Range("I122").Value = 2400.2856
$ws.Range("K122").Value = 7200.8568
$ws.Range("M122").Value = -4750.8568
$ws.Range("H132").Value = 6837.077
$ws.Range("I132").Value = 5220.4
$ws.Range("J132").Value = 7847.5
$ws.Range("K132").Value = 15661.2
$ws.Range("L132").Value = 23542.5
$ws.Range("M132").Value = -13131.2
$ws.Range("N132").Value = -28602.5
$ws.Range("H136").Value = 3795.6206
$ws.Range("J136").Value = 7052.6924
$ws.Range("L136").Value = 21158.0772
$ws.Range("N136").Value = -26258.0772
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1524528.1
$ws.Range("I126").Value = 3833.3333
$ws.Range("K126").Value = 11499.9999
$ws.Range("M126").Value = -9029.999899999999
$ws.Range("H132").Value = 7250446.5
$ws.Range("I132").Value = 4874.346
$ws.Range("K132").Value = 14623.038
$ws.Range("M132").Value = -12093.038
$ws.Range("H136").Value = 8343.083000000001
$ws.Range("I136").Value = 7934.933
$ws.Range("J136").Value = 9023.333000000001
$ws.Range("K136").Value = 23804.799
$ws.Range("L136").Value = 27069.999
$ws.Range("M136").Value = -21254.799
$ws.Range("N136").Value = -32169.999
